$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# The 6-row value/style block that gets appended (B2:B7 pattern).
$srcBlock = $ws1.Range("B2:B7")

# Append a copy of the B2:B7 block to Sheet1 starting at row 11.
$srcBlock.Copy()
$dst1 = $ws1.Range("B11:B16")
$dst1.PasteSpecial(-4163)   # xlPasteValues
$dst1.PasteSpecial(-4122)   # xlPasteFormats

# Create Sheet2 as a copy of Sheet1 (placed right after it), then wipe its
# contents so it inherits the same sheet formatting/phonetic settings.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"
$ws2.Cells.Clear()

# Fill Sheet2 with the same value/style block, starting at row 3.
$srcBlock.Copy()
$dst2 = $ws2.Range("B3:B8")
$dst2.PasteSpecial(-4163)   # xlPasteValues
$dst2.PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("B3:B8").Select()

# Leave Sheet1 active with B9 selected, matching the final UI state.
$ws1.Activate()
$ws1.Range("B9").Select()
